$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old row 3 (copy_processo_2 / 4835245-15.2024.8.01.2832 / ...) becomes the
# new row 4, shifted down one row (with its "Arquivo" label simplified from
# "copy_processo_2" to "processo_2"). Copy the whole row (values + formats)
# down to row 4 first so styles / string-vs-number cell types are preserved.
$ws.Range("A3:F3").Copy($ws.Range("A4:F4"))

# Fix up the "Arquivo" label on the row that moved down to row 4.
$ws.Range("A4").Value = "processo_2"

# Now overwrite row 3 with the new record that was inserted ahead of it.
$ws.Range("A3").Value = "copy_processo_3"
$ws.Range("B3").Value = "3130687-11.2024.8.01.5042"
$ws.Range("C3").Value = "Nome Aleatório 98"
$ws.Range("D3").Value = "Advogado Exemplo"
$ws.Range("E3").Value = 43679
$ws.Range("F3").Value = "7/8/2024"
